# Generate Report for Handback
# Replace the old UUID-based file identifiers/timestamps with the new ones
# across the Overview, zh-cn, and de-de worksheets.

$wb = $excel.ActiveWorkbook

$oldId1 = "674e4054-463e-4f57-b03f-a75f83fdb0ff"
$newId1 = "8022f275-f833-4477-b9d8-4ea2be31bf2d"
$oldId2 = "c57c4e29-70fe-47d4-9162-3e3b2459bc17"
$newId2 = "ffff5fe05572-33a6-44b7-b8b8-40f801a4d732"

$oldHash1 = "fc4e56ee4a7d451e7532a349b2707d279729f771"
$newHash1 = "816eed1bbd4b688c757dbb3800454b802dd95472"

# -------------------- Overview sheet --------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId1.md"
$wsOverview.Range("B2").Value = "e2e\$newId1.md"
$wsOverview.Range("G2").Value = "2016-09-01 13:10:36"

$wsOverview.Range("A3").Value = "$newId2.md"
$wsOverview.Range("B3").Value = "e2e\$newId2.md"
$wsOverview.Range("G3").Value = "2016-09-01 13:10:36"

# -------------------- zh-cn sheet --------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId1.md"
$wsZhCn.Range("I2").Value = "$newId1.md"
$wsZhCn.Range("G2").Value = "$newId1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("J2").Value = "$newId1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-01 13:10:30"
$wsZhCn.Range("K2").Value = "2016-09-01 13:10:59"

$wsZhCn.Range("A3").Value = "$newId2.md"
$wsZhCn.Range("I3").Value = "$newId2.md"
$wsZhCn.Range("G3").Value = "$newId1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("J3").Value = "$newId1.$newHash1.zh-cn.xlf"
$wsZhCn.Range("H3").Value = "2016-09-01 13:10:30"
$wsZhCn.Range("K3").Value = "2016-09-01 13:10:59"

# -------------------- de-de sheet --------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId1.md"
$wsDeDe.Range("I2").Value = "$newId1.md"
$wsDeDe.Range("G2").Value = "$newId1.$newHash1.de-de.xlf"
$wsDeDe.Range("J2").Value = "$newId1.$newHash1.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-01 13:10:36"
$wsDeDe.Range("K2").Value = "2016-09-01 13:11:22"

$wsDeDe.Range("A3").Value = "$newId2.md"
$wsDeDe.Range("I3").Value = "$newId2.md"
$wsDeDe.Range("G3").Value = "$newId1.$newHash1.de-de.xlf"
$wsDeDe.Range("J3").Value = "$newId1.$newHash1.de-de.xlf"
$wsDeDe.Range("H3").Value = "2016-09-01 13:10:36"
$wsDeDe.Range("K3").Value = "2016-09-01 13:11:22"

# -------------------- Hyperlinks --------------------
# Overview sheet hyperlinks (B2, B3)
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newId1.md"
$wsOverview.Hyperlinks.Item(2).TextToDisplay = "e2e\$newId2.md"

# zh-cn sheet hyperlinks (A2, I2, A3, I3)
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$wsZhCn.Hyperlinks.Item(2).TextToDisplay = "$newId1.md"
$wsZhCn.Hyperlinks.Item(3).TextToDisplay = "$newId2.md"
$wsZhCn.Hyperlinks.Item(4).TextToDisplay = "$newId2.md"

# de-de sheet hyperlinks (A2, I2, A3, I3)
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newId1.md"
$wsDeDe.Hyperlinks.Item(2).TextToDisplay = "$newId1.md"
$wsDeDe.Hyperlinks.Item(3).TextToDisplay = "$newId2.md"
$wsDeDe.Hyperlinks.Item(4).TextToDisplay = "$newId2.md"
